$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 header cells Z9:AI9 -- copy text-typed numbering labels from M9:V9
$ws.Range("M9:V9").Copy($ws.Range("Z9"))

# Row 10
$ws.Range("Y10").Value = @'
SELECT sum(mitarbeiteranzahl) AS MA_Bundesland
 FROM standorte
WHERE bundesland ='Berlin'

'@
$ws.Range("Y10").WrapText = $true
$ws.Range("Z10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2984 ms, verstrichene Zeit = 382 ms.
'@
$ws.Range("Z10").WrapText = $true
$ws.Range("AA10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2937 ms, verstrichene Zeit = 375 ms.
'@
$ws.Range("AA10").WrapText = $true
$ws.Range("AB10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2969 ms, verstrichene Zeit = 374 ms.
'@
$ws.Range("AB10").WrapText = $true
$ws.Range("AC10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2921 ms, verstrichene Zeit = 376 ms.

'@
$ws.Range("AC10").WrapText = $true
$ws.Range("AD10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2937 ms, verstrichene Zeit = 380 ms.
'@
$ws.Range("AD10").WrapText = $true
$ws.Range("AE10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2937 ms, verstrichene Zeit = 373 ms.
'@
$ws.Range("AE10").WrapText = $true
$ws.Range("AF10").Value = @'

 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2968 ms, verstrichene Zeit = 374 ms.
'@
$ws.Range("AF10").WrapText = $true
$ws.Range("AG10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2922 ms, verstrichene Zeit = 375 ms.
'@
$ws.Range("AG10").WrapText = $true
$ws.Range("AH10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2953 ms, verstrichene Zeit = 374 ms.
'@
$ws.Range("AH10").WrapText = $true
$ws.Range("AI10").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2985 ms, verstrichene Zeit = 374 ms.
'@
$ws.Range("AI10").WrapText = $true

# Row 11
$ws.Range("Z11").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2484 ms, verstrichene Zeit = 324 ms.
'@
$ws.Range("Z11").WrapText = $true
$ws.Range("Y11").Value = @'
SELECT min(betraggesamt)AS Kleinster_Gesamt_Betrag 
 FROM bestellung
WHERE menge <= 100 

'@
$ws.Range("Y11").WrapText = $true
$ws.Range("AA11").Value = @'
SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2578 ms, verstrichene Zeit = 322 ms.
'@
$ws.Range("AA11").WrapText = $true
$ws.Range("AB11").Value = @'

 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2592 ms, verstrichene Zeit = 319 ms.
'@
$ws.Range("AB11").WrapText = $true
$ws.Range("AC11").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2454 ms, verstrichene Zeit = 319 ms.

'@
$ws.Range("AC11").WrapText = $true
$ws.Range("AD11").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2578 ms, verstrichene Zeit = 319 ms.
'@
$ws.Range("AD11").WrapText = $true
$ws.Range("AE11").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2484 ms, verstrichene Zeit = 321 ms.
'@
$ws.Range("AE11").WrapText = $true
$ws.Range("AF11").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2579 ms, verstrichene Zeit = 319 ms.
'@
$ws.Range("AF11").WrapText = $true
$ws.Range("AG11").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2483 ms, verstrichene Zeit = 319 ms.
'@
$ws.Range("AG11").WrapText = $true
$ws.Range("AH11").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2471 ms, verstrichene Zeit = 319 ms.
'@
$ws.Range("AH11").WrapText = $true
$ws.Range("AI11").Value = @'
SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2500 ms, verstrichene Zeit = 331 ms.
'@
$ws.Range("AI11").WrapText = $true

# Row 12
$ws.Range("AH12").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2969 ms, verstrichene Zeit = 374 ms.
'@
$ws.Range("AH12").WrapText = $true
$ws.Range("Y12").Value = @'
SELECT max(preis),min(preis)
 FROM bestellung
WHERE menge <= 100
'@
$ws.Range("Y12").WrapText = $true
$ws.Range("Z12").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2969 ms, verstrichene Zeit = 369 ms.
'@
$ws.Range("Z12").WrapText = $true
$ws.Range("AA12").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2844 ms, verstrichene Zeit = 373 ms.

'@
$ws.Range("AA12").WrapText = $true
$ws.Range("AB12").Value = @'

 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2937 ms, verstrichene Zeit = 372 ms.
'@
$ws.Range("AB12").WrapText = $true
$ws.Range("AC12").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2858 ms, verstrichene Zeit = 372 ms.
'@
$ws.Range("AC12").WrapText = $true
$ws.Range("AD12").Value = @'

 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2969 ms, verstrichene Zeit = 370 ms.
'@
$ws.Range("AD12").WrapText = $true
$ws.Range("AE12").Value = @'
SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2952 ms, verstrichene Zeit = 372 ms.
'@
$ws.Range("AE12").WrapText = $true
$ws.Range("AF12").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2970 ms, verstrichene Zeit = 370 ms.
'@
$ws.Range("AF12").WrapText = $true
$ws.Range("AG12").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3015 ms, verstrichene Zeit = 380 ms.
'@
$ws.Range("AG12").WrapText = $true
$ws.Range("AI12").Value = @'
SQL Server-Ausführungszeiten: 
, CPU-Zeit = 2923 ms, verstrichene Zeit = 374 ms.
'@
$ws.Range("AI12").WrapText = $true

# Row 13
$ws.Range("Y13").Value = @'
SELECT max(fahrzeuganzahl),max(mitarbeiteranzahl) 
 FROM lieferdienst
WHERE fahrzeugtyp = 'Auto'

'@
$ws.Range("Y13").WrapText = $true
$ws.Range("Z13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3202 ms, verstrichene Zeit = 408 ms.
'@
$ws.Range("Z13").WrapText = $true
$ws.Range("AA13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3204 ms, verstrichene Zeit = 404 ms.

'@
$ws.Range("AA13").WrapText = $true
$ws.Range("AB13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3094 ms, verstrichene Zeit = 404 ms.
'@
$ws.Range("AB13").WrapText = $true
$ws.Range("AC13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3172 ms, verstrichene Zeit = 407 ms.
'@
$ws.Range("AC13").WrapText = $true
$ws.Range("AD13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3186 ms, verstrichene Zeit = 404 ms.
'@
$ws.Range("AD13").WrapText = $true
$ws.Range("AE13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3170 ms, verstrichene Zeit = 404 ms.
'@
$ws.Range("AE13").WrapText = $true
$ws.Range("AF13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3048 ms, verstrichene Zeit = 404 ms.
'@
$ws.Range("AF13").WrapText = $true
$ws.Range("AG13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3203 ms, verstrichene Zeit = 407 ms.
'@
$ws.Range("AG13").WrapText = $true
$ws.Range("AH13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3172 ms, verstrichene Zeit = 405 ms.
'@
$ws.Range("AH13").WrapText = $true
$ws.Range("AI13").Value = @'
 SQL Server-Ausführungszeiten: 
, CPU-Zeit = 3202 ms, verstrichene Zeit = 405 ms.
'@
$ws.Range("AI13").WrapText = $true

# Update selection to match the final active cell
$ws.Range("AI14").Select()
